$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Change 1: " eschaufer avecq la flamme de la " -> split so the run
# becomes " eschaufe" + "s" + " avecq la flamme de la ", where the new
# "s" run carries no w:color override (only inherited formatting, i.e.
# rtl=0), unlike its colored ("000000") neighbours.
#
# The object model always stamps an explicit color when Font.Color is
# set (even wdColorAutomatic -> w:color val="auto"), so instead we
# borrow an existing plain run's formatting (a lone "q" in "quarres"
# elsewhere in the document, whose rPr is just <w:rtl w:val="0"/>,
# no color) via Copy/Paste, which preserves that plain formatting,
# and then fix up the pasted character's text to "s".
# ------------------------------------------------------------------

$srcFind = $d.Content
$srcFind.Find.Execute("quarres", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$srcChar = $d.Range($srcFind.Start, $srcFind.Start + 1)
$srcChar.Copy() | Out-Null

$tgt = $d.Content
$tgt.Find.Execute("eschaufer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailStart = $tgt.End - 1
$tailEnd = $tgt.End

# Drop the trailing "r" of "eschaufer", leaving " eschaufe".
$rTail = $d.Range($tailStart, $tailEnd)
$rTail.Text = ""

# Paste the plain-formatted character right after " eschaufe".
$insertHere = $d.Range($tailStart, $tailStart)
$insertHere.Paste() | Out-Null

# Turn the pasted "q" into the "s" we actually need, keeping its
# (color-less) run formatting intact.
$pasted = $d.Range($tailStart, $tailStart + 1)
$pasted.Text = "s"

# ------------------------------------------------------------------
# Change 2: "seicher bien alaise Apres fortifie les couches" ->
# "seicher bien a laise Apres fortifie les couches" (space inserted
# between "a" and "laise").
# ------------------------------------------------------------------

$d.Content.Find.Execute("seicher bien alaise Apres fortifie les couches", $true, $false, $false, $false, $false, $true, 1, $false, "seicher bien a laise Apres fortifie les couches", 2)
